$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "From the Heart Tutoring" bullet - expand the launch sentence with
# a new lead-in clause and spell out "JavaScript".
#   Paragraph 25 (1-based Word paragraph index):
#   "Launched the company website using HTML5, CSS, and JS on AWS S3,
#    cutting onboarding time by 99%."
#   ->
#   "Lead all software development and launched the company website using
#    HTML5, CSS, and JavaScript on AWS S3, cutting onboarding time by 99%."
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(25)
$old1 = "Launched the company website using HTML5, CSS, and JS on AWS S3, cutting onboarding time by 99%."
$new1 = "Lead all software development and launched the company website using HTML5, CSS, and JavaScript on AWS S3, cutting onboarding time by 99%."
$p1.Range.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: replace the "Designed and deployed a MVC booking service..."
# bullet entirely with a tutoring-session bullet, and add w:after="94" to the
# paragraph spacing (4.7pt, since spacing is stored in twentieths of a point).
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(26)
$old2 = "Designed and deployed a MVC booking service with a React frontend, C# REST API backend, and DynamoDB data layer, containerized with Docker and hosted using AWS microservices"
$new2 = "Delivered 1,100+ tutoring sessions in engineering, computer science, and STEM subjects."
$p2.Range.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null
$p2.SpaceAfter = 4.7

# ---------------------------------------------------------------------------
# Change 3: project title "AI Drone Detection System" ->
# "Full-Stack Job Application Assistant"
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(41)
$p3.Range.Find.Execute("AI Drone Detection System", $false, $false, $false, $false, $false, $true, 1, $false, "Full-Stack Job Application Assistant", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 4: project subtitle "Texas State University" -> "Personal"
# (scoped to its own paragraph so the other 4 occurrences of this phrase
# elsewhere in the resume are left untouched)
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(42)
$p4.Range.Find.Execute("Texas State University", $false, $false, $false, $false, $false, $true, 1, $false, "Personal", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 5 & 6: rewrite the drone-detection bullet to describe the job
# application assistant project instead.
#   "Developed six Python-based models using convolutional neural networks
#    (CNNs) and random forests to classify audio clips from three datasets
#    as drone or non-drone sourced."
#   ->
#   "Developed and deployed a full-stack application using React, Tailwind
#    CSS, C#, and SQL that leverages AI to assist software developers with
#    job applications."
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(43)
$p5.Range.Find.Execute("six ", $false, $false, $false, $false, $false, $true, 1, $false, "and deployed ", 2) | Out-Null

$p6 = $d.Paragraphs.Item(43)
$old6 = "Python-based models using convolutional neural networks (CNNs) and random forests to classify audio clips from three datasets as drone or non-drone sourced"
$new6 = "a full-stack application using React, Tailwind CSS, C#, and SQL that leverages AI to assist software developers with job applications"
$p6.Range.Find.Execute($old6, $false, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null
